$wb = $excel.ActiveWorkbook

# --- Sheet "Means" ---
$wsMeans = $wb.Worksheets.Item("Means")

# Row 9: Total Cancer Risk (per million)
$wsMeans.Range("B9").Value = 26
$wsMeans.Range("C9").Value = 27
$wsMeans.Range("D9").Value = 29

# Row 10: Total Respiratory (hazard quotient)
$wsMeans.Range("B10").Value = 0.31
$wsMeans.Range("C10").Value = 0.34
$wsMeans.Range("D10").Value = 0.39
$wsMeans.Range("E10").Value = 0.35
$wsMeans.Range("F10").Value = 0.34
$wsMeans.Range("G10").Value = 0.33

# --- Sheet "Standard Deviations" ---
$wsSD = $wb.Worksheets.Item("Standard Deviations")

# Row 9: Total Cancer Risk (per million)
$wsSD.Range("B9").Value = 8.3
$wsSD.Range("C9").Value = 7.2
$wsSD.Range("D9").Value = 4.6
$wsSD.Range("E9").Value = 1.5
$wsSD.Range("F9").Value = 1.1
$wsSD.Range("G9").Value = 1.6

# Row 10: Total Respiratory (hazard quotient)
$wsSD.Range("B10").Value = 0.11
$wsSD.Range("C10").Value = 0.1
$wsSD.Range("E10").Value = 0.05
$wsSD.Range("F10").Value = 0.05
$wsSD.Range("G10").Value = 0.048
